# Update "想去人数" (interested-count) figures pulled from a fresh data export.
$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 3419
$ws1.Range("F3").Value = 743
$ws1.Range("F5").Value = 6973
$ws1.Range("F6").Value = 2447
$ws1.Range("F10").Value = 38
$ws1.Range("F11").Value = 77
$ws1.Range("F14").Value = 571

# Sheet "全部类型" (all types) mirrors the same events at shifted row numbers
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 3419
$ws4.Range("F4").Value = 743
$ws4.Range("F6").Value = 6973
$ws4.Range("F7").Value = 2447
$ws4.Range("F11").Value = 38
$ws4.Range("F12").Value = 77
$ws4.Range("F15").Value = 571
